# Add more tests; clean more code
$wb = $excel.ActiveWorkbook

$wsTrend = $wb.Worksheets.Item("Trend_instructions")

# Header B1: "agg_fuel" -> "CEDS_fuel"
$wsTrend.Range("B1").Value = "CEDS_fuel"

# Data B2: "coal" -> "coal coke"
$wsTrend.Range("B2").Value = "coal coke"

# Update the active selection on the Trend_instructions sheet to B2
$wsTrend.Activate()
$wsTrend.Range("B2").Select()
